# Update the "nao_enviados" sheet:
#  - row 8's Telefone1 (B8) was stored as text; turn it into a real number,
#    matching the other numeric phone cells above it.
#  - append the rest of the captured "nao enviados" log (rows 9-37) coming
#    from the improved disparo.py run (separate image/PDF dispatch).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing row 8 --------------------------------------------------
$ws.Cells.Item(8, 2).Value = 5561982757272

# --- append new rows 9-37 -------------------------------------------------
$rows = @(
    @{ Row = 9;  Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 10; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 11; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 12; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 13; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 14; Name = "thiago ";   Phone = "619187913" },
    @{ Row = 15; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 16; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 17; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 18; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 19; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 20; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 21; Name = "thiago ";   Phone = "619187913" },
    @{ Row = 22; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 23; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 24; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 25; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 26; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 27; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 28; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 29; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 30; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 31; Name = "thiago ";   Phone = "619187913" },
    @{ Row = 32; Name = "Ana";       Phone = "(61)98182-6392" },
    @{ Row = 33; Name = "dhiogenes"; Phone = "5561982757272" },
    @{ Row = 34; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 35; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 36; Name = "guilherme"; Phone = "619850276" },
    @{ Row = 37; Name = "thiago ";   Phone = "619187913" }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Name (column A) - always plain text.
    $ws.Cells.Item($r, 1).Value = "'" + $item.Name
    $ws.Cells.Item($r, 1).Style = "Normal"

    # Phone (column B) - kept as text even when it looks numeric, matching
    # the source log (only B8 above was normalised to a real number).
    $ws.Cells.Item($r, 2).Value = "'" + $item.Phone
    $ws.Cells.Item($r, 2).Style = "Normal"

    # Telefone2 (column C) stays blank, like the rest of the sheet.
    $ws.Cells.Item($r, 3).Value = "'"
    $ws.Cells.Item($r, 3).Style = "Normal"
}
